$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 07:30"

# Row 15 - Pakistan
$ws.Range("B15").Value = 213470
$ws.Range("C15").Value = 4133
$ws.Range("D15").Value = 100802
$ws.Range("E15").Value = 108273
$ws.Range("G15").Value = 91
$ws.Range("H15").Value = 4395

# Row 17 - Alemania
$ws.Range("D17").Value = 179800
$ws.Range("E17").Value = 6980

# Row 74 - Uzbekistan
$ws.Range("B74").Value = 8627
$ws.Range("C74").Value = 124
$ws.Range("E74").Value = 2919

# Row 86 - Kirguistan
$ws.Range("B86").Value = 5735
$ws.Range("C86").Value = 439
$ws.Range("E86").Value = 3230
$ws.Range("G86").Value = 5
$ws.Range("H86").Value = 62

# Row 99 - Tailandia
$ws.Range("B99").Value = 3173
$ws.Range("C99").Value = 2
$ws.Range("D99").Value = 3059
$ws.Range("E99").Value = 56
